$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data record was added for "Berenjena" (eggplant) prices.
# It belongs right before the existing row 168, so shift rows 168:197
# down by one (to 169:198) and populate the freed-up row 168 with the
# new record's values.
$ws.Rows.Item(168).Insert()

$ws.Cells.Item(168, 1).Value = 10
$ws.Cells.Item(168, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(168, 3).Value = "La Araucanía"
$ws.Cells.Item(168, 4).Value = 44504
$ws.Cells.Item(168, 5).Value = 9
$ws.Cells.Item(168, 6).Value = 100112001
$ws.Cells.Item(168, 7).Value = "Berenjena"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 75
$ws.Cells.Item(168, 11).Value = 10000
$ws.Cells.Item(168, 12).Value = 12000
$ws.Cells.Item(168, 13).Value = 10933
$ws.Cells.Item(168, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(168, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(168, 16).Value = 182
$ws.Cells.Item(168, 17).Value = 60
$ws.Cells.Item(168, 18).Value = "Hortaliza"
